$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I0 (I) and IF (J)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-18
$data = @(
    @(8, 9),
    @(8, 8),
    @(6, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(4, 4),
    @(7, 7),
    @(6, 6),
    @(2, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
